$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "Anselmo-Gestão Integr"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "André Guimarães-Desenho Técn"

# Row 3
$ws.Range("B3").Value = "[-, Valmir-Metrologia, -, -]"
$ws.Range("C3").Value = "Anselmo-Gestão Integr"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "André Guimarães-Desenho Técn"

# Row 4
$ws.Range("B4").Value = "[-, Valmir-Metrologia, -, -]"
$ws.Range("C4").Value = "[-, Joel Lima-Tecnologia dos Materiais]"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "-"

# Row 6
$ws.Range("B6").Value = "[-, Valmir-Metrologia, -, -]"
$ws.Range("C6").Value = "[-, Joel Lima-Tecnologia dos Materiais]"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "-"

# Row 7
$ws.Range("C7").Value = "[-, Joel Lima-Tecnologia dos Materiais]"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "Cleidson-Circuitos Elétrico"

# Row 8
$ws.Range("C8").Value = "[-, Joel Lima-Tecnologia dos Materiais]"
$ws.Range("E8").Value = "[Valmir-Metrologia, -, -, -]"
$ws.Range("F8").Value = "Cleidson-Circuitos Elétrico"
